# "Added week 13 Info"
#
# Week 13 previously had only the header row + the totals rows (20/21).
# Three new time entries are added (rows 2-4), which changes the Weekly
# Total (E20) and therefore cascades through the Project Total (E21) of
# every later week (14, 15) and the Final sheet, since each week's E21
# is "=E20+'<previous week>'!E21".
#
# Week 12 was the previously-active/selected sheet; Week 13 becomes the
# new active sheet with E3 selected.

$wb = $excel.ActiveWorkbook

$weprevious = $wb.Worksheets.Item("Week 12")
$week13     = $wb.Worksheets.Item("Week 13")
$final      = $wb.Worksheets.Item("Final")

# Bring over the existing date/time number formatting (styles 5 & 6) and
# the description-column formatting (style 7) by copying a same-shaped
# block of already-styled cells from Week 12, then overwrite the copied
# values below. This reuses the workbook's existing style indices
# instead of minting new ones.
$weprevious_row2to3 = $weprevious.Range("A2:E3")
$weprevious_row2to3.Copy($week13.Range("A2"))
$weprevious.Range("A2:E2").Copy($week13.Range("A4"))
$excel.CutCopyMode = 0

# Row 2 - 4/6 (Saturday)
$week13.Range("A2").Value = 43560
$week13.Range("B2").Value = 0.52083333333333337
$week13.Range("C2").Value = 0.58333333333333337
$week13.Range("D2").Value = "Worked on Wish List feature"
$week13.Range("E2").Value = 1.5

# Row 4 - 4/13 (set before row 3 so new shared strings are written in
# the same first-seen order as the source workbook: "...Wish List
# feature", "Finished Integrating Search Feature", "...Wish List and
# Search Feature")
$week13.Range("A4").Value = 43567
$week13.Range("B4").Value = 0.75
$week13.Range("C4").Value = 0.83333333333333337
$week13.Range("D4").Value = "Finished Integrating Search Feature"
$week13.Range("E4").Value = 2

# Row 3 - 4/11
$week13.Range("A3").Value = 43565
$week13.Range("B3").Value = 0.52083333333333337
$week13.Range("C3").Value = 0.63541666666666663
$week13.Range("D3").Value = "Worked on Wish List and Search Feature"
$week13.Range("E3").Value = 2.75

# Nudge the Final sheet's scroll position without disturbing its
# current selection (activeCell stays A2).
$final.Activate()
$finalWindow = $excel.ActiveWindow
$finalWindow.ScrollRow = 5
$finalWindow.ScrollColumn = 1
$final.Range("A2").Select()

# Week 13 becomes the selected/active tab (Week 12 loses tabSelected
# automatically since only one sheet can be active at a time), with
# E3 as the active cell.
$week13.Activate()
$week13.Range("E3").Select()
